$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last three lineage rows (full_name, last_name, orders) collapse away;
# "orders" now lands on what used to be the first_name row (row 30), and the
# old rows 31-33 are removed entirely.
$ws.Range("D30").Value = "orders"
$ws.Range("E30").Value = "orders"
$ws.Rows("31:33").Delete()
